$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15: Morning Glass of Ether | Ether
$ws.Range("H15").Value = 4762.8096
$ws.Range("I15").Value = 4762.8096
$ws.Range("K15").Value = 14288.4288
$ws.Range("M15").Value = -14119.4288

# Row 125: Body over Mind | Grade 5 Dexterity Alkahest
$ws.Range("H125").Value = 1658.5714
$ws.Range("J125").Value = 1736.5
$ws.Range("L125").Value = 15628.5
$ws.Range("N125").Value = -20548.5

# Row 129: Practical Command | Commanding Craftsman's Draught
$ws.Range("H129").Value = 19906.547
$ws.Range("I129").Value = 1326.3636
$ws.Range("J129").Value = 24772.785
$ws.Range("K129").Value = 3979.0908
$ws.Range("L129").Value = 74318.355
$ws.Range("M129").Value = 1020.9092
$ws.Range("N129").Value = -84318.355

# Row 132: Fast-forwarding Flora | Growth Formula Lambda
$ws.Range("H132").Value = 45521.59
$ws.Range("I132").Value = 34676.844
$ws.Range("J132").Value = 74440.914
$ws.Range("K132").Value = 104030.532
$ws.Range("L132").Value = 223322.742
$ws.Range("M132").Value = -101500.532
$ws.Range("N132").Value = -228382.742

# Row 137: Cutting Edge of Culinary Quality | Magnesia Whetstone
$ws.Range("H137").Value = 3216673.2
$ws.Range("I137").Value = 8563186
$ws.Range("K137").Value = 25689558
$ws.Range("M137").Value = -25687008

$ws = $wb.Worksheets.Item("ARM")
# Row 61: Dealing with the Tough Stuff | Cobalt Ingot
$ws.Range("H61").Value = 1982.8214
$ws.Range("I61").Value = 1207.3684
$ws.Range("J61").Value = 3619.889
$ws.Range("K61").Value = 1207.3684
$ws.Range("L61").Value = 3619.889
$ws.Range("M61").Value = -995.3684000000001
$ws.Range("N61").Value = -4043.889

# Row 132: Don't Bore Me, Ore Me | Mountain Chromite Ingot
$ws.Range("H132").Value = 12196904
$ws.Range("I132").Value = 15152911
$ws.Range("K132").Value = 45458733
$ws.Range("M132").Value = -45456203

# Row 136: Metal with Mettle | Cobalt Tungsten Ingot
$ws.Range("H136").Value = 1982.8214
$ws.Range("I136").Value = 1207.3684
$ws.Range("J136").Value = 3619.889
$ws.Range("K136").Value = 3622.1052
$ws.Range("L136").Value = 10859.667
$ws.Range("M136").Value = -1072.1052
$ws.Range("N136").Value = -15959.667

$ws = $wb.Worksheets.Item("BSM")
# Row 105: Ingot to Wing It | Molybdenum Ingot
$ws.Range("H105").Value = 2400.875
$ws.Range("I105").Value = 2287.3914
$ws.Range("J105").Value = 5011
$ws.Range("K105").Value = 2287.3914
$ws.Range("L105").Value = 5011
$ws.Range("M105").Value = -540.3914
$ws.Range("N105").Value = -8505

# Row 134: Ruthenium Supremium | Ruthenium Ingot
$ws.Range("H134").Value = 2186.4055
$ws.Range("I134").Value = 1620
$ws.Range("K134").Value = 4860
$ws.Range("M134").Value = -2325

$ws = $wb.Worksheets.Item("CRP")
# Row 58: You Do the Heavy Lifting | Mahogany Lumber
$ws.Range("H58").Value = 2733.9583
$ws.Range("I58").Value = 1600.9231
$ws.Range("J58").Value = 4073
$ws.Range("K58").Value = 1600.9231
$ws.Range("L58").Value = 4073
$ws.Range("M58").Value = -1397.9231
$ws.Range("N58").Value = -4479

# Row 136: Turali Quality | Dark Mahogany Lumber
$ws.Range("H136").Value = 2733.9583
$ws.Range("I136").Value = 1600.9231
$ws.Range("J136").Value = 4073
$ws.Range("K136").Value = 4802.7693
$ws.Range("L136").Value = 12219
$ws.Range("M136").Value = -2252.7693
$ws.Range("N136").Value = -17319

$ws = $wb.Worksheets.Item("CUL")
# Row 5: What a Sap | Maple Syrup
$ws.Range("H5").Value = 8908.666999999999
$ws.Range("I5").Value = 12838
$ws.Range("J5").Value = 1050
$ws.Range("K5").Value = 38514
$ws.Range("L5").Value = 3150
$ws.Range("M5").Value = -38402
$ws.Range("N5").Value = -3374

# Row 34: Fever Pitch | Chamomile Tea
$ws.Range("H34").Value = 1527.1765
$ws.Range("I34").Value = 462.5
$ws.Range("J34").Value = 1854.7693
$ws.Range("K34").Value = 1387.5
$ws.Range("L34").Value = 5564.3079
$ws.Range("M34").Value = -1303.5
$ws.Range("N34").Value = -5732.3079

# Row 131: The Mountain Steeped | Tsai tou Vounou
$ws.Range("H131").Value = 111530.6
$ws.Range("I131").Value = 150513.62
$ws.Range("J131").Value = 85541.914
$ws.Range("K131").Value = 451540.86
$ws.Range("L131").Value = 256625.742
$ws.Range("M131").Value = -446500.86
$ws.Range("N131").Value = -266705.742

# Row 135: Not-so-secret Ingredient | Royal Maple Syrup
$ws.Range("H135").Value = 8908.666999999999
$ws.Range("I135").Value = 12838
$ws.Range("J135").Value = 1050
$ws.Range("K135").Value = 115542
$ws.Range("L135").Value = 9450
$ws.Range("M135").Value = -113007
$ws.Range("N135").Value = -14520

# Row 137: Creative Chocolate | Gateau au Chocolat
$ws.Range("H137").Value = 9282.963
$ws.Range("I137").Value = 4357.778
$ws.Range("J137").Value = 11745.556
$ws.Range("K137").Value = 13073.334
$ws.Range("L137").Value = 35236.66800000001
$ws.Range("M137").Value = -7973.334000000001
$ws.Range("N137").Value = -45436.66800000001

$ws = $wb.Worksheets.Item("GSM")
# Row 97: If I'd a Koppranickel for Every Time... | Koppranickel Ingot
$ws.Range("H97").Value = 3601.0476
$ws.Range("I97").Value = 2482.353
$ws.Range("K97").Value = 2482.353
$ws.Range("M97").Value = -1986.353

# Row 102: Put the Metal to the Peddle | Durium Ingot
$ws.Range("H102").Value = 971.8461
$ws.Range("I102").Value = 862
$ws.Range("J102").Value = 1338
$ws.Range("K102").Value = 862
$ws.Range("L102").Value = 1338
$ws.Range("M102").Value = 760
$ws.Range("N102").Value = -4582

# Row 126: Gold Rush Order | Phrygian Gold Ingot
$ws.Range("H126").Value = 3127.8462
$ws.Range("I126").Value = 3364
$ws.Range("J126").Value = 2925.4285
$ws.Range("K126").Value = 10092
$ws.Range("L126").Value = 8776.2855
$ws.Range("M126").Value = -7622
$ws.Range("N126").Value = -13716.2855

# Row 132: On Board for Lar | Lar Ingot
$ws.Range("H132").Value = 27029430
$ws.Range("I132").Value = 41668216
$ws.Range("K132").Value = 125004648
$ws.Range("M132").Value = -125002118

$ws = $wb.Worksheets.Item("LTW")
# Row 22: Skin off Their Backs | Aldgoat Leather
$ws.Range("H22").Value = 998
$ws.Range("I22").Value = 1025.7142
$ws.Range("J22").Value = 933.3333
$ws.Range("K22").Value = 1025.7142
$ws.Range("L22").Value = 933.3333
$ws.Range("M22").Value = -730.7141999999999
$ws.Range("N22").Value = -1523.3333

# Row 27: Fire and Hide | Aldgoat Leather
$ws.Range("H27").Value = 998
$ws.Range("I27").Value = 1025.7142
$ws.Range("J27").Value = 933.3333
$ws.Range("K27").Value = 1025.7142
$ws.Range("L27").Value = 933.3333
$ws.Range("M27").Value = -918.7141999999999
$ws.Range("N27").Value = -1147.3333

# Row 100: Tiger in the Sack | Tiger Leather
$ws.Range("H100").Value = 3890.158
$ws.Range("I100").Value = 2950.9
$ws.Range("J100").Value = 4933.778
$ws.Range("K100").Value = 2950.9
$ws.Range("L100").Value = 4933.778
$ws.Range("M100").Value = -2409.9
$ws.Range("N100").Value = -6015.778

# Row 119: Fit for a Friend | Swallowskin Gloves of Fending
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()

# Row 132: Tenets of Tanning | Silver Lobo Leather
$ws.Range("H132").Value = 3659.0527
$ws.Range("I132").Value = 3127.65
$ws.Range("J132").Value = 4249.5
$ws.Range("K132").Value = 9382.950000000001
$ws.Range("L132").Value = 12748.5
$ws.Range("M132").Value = -6852.950000000001
$ws.Range("N132").Value = -17808.5

$ws = $wb.Worksheets.Item("WVR")
# Row 70: An Account of My Boots | Holy Rainbow Shoes
$ws.Range("H70").Value = 34952.5
$ws.Range("J70").Value = 34952.5
$ws.Range("L70").Value = 34952.5
$ws.Range("N70").Value = -35582.5

# Row 73: Soot in My Hair and Scars on My Feet (L) | Holy Rainbow Shoes
$ws.Range("H73").Value = 34952.5
$ws.Range("J73").Value = 34952.5
$ws.Range("L73").Value = 34952.5
$ws.Range("N73").Value = -37136.5

# Row 132: Comfy Cabins | Snow Cotton Cloth
$ws.Range("H132").Value = 1116790.2
$ws.Range("I132").Value = 1554233.5
$ws.Range("J132").Value = 3298.2727
$ws.Range("K132").Value = 4662700.5
$ws.Range("L132").Value = 9894.8181
$ws.Range("M132").Value = -4660170.5
$ws.Range("N132").Value = -14954.8181
